# Update column G ("K") values on Sheet1 per the commit:
# "regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2  = 2
    3  = 0
    4  = 0
    6  = 0
    7  = 1
    8  = 1
    9  = 1
    11 = 0
    12 = 1
    13 = 1
    14 = 0
    15 = 2
    17 = 2
    18 = 1
    19 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}

$wb.Save()
